$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The RB006/RB007/RB008 test-batch codes in column C were consolidated to the
# already-used "RB005" code (commit: "modified test data and rm 'N.A' in
# configuration") -- this also makes the three now-unused shared strings drop
# out of the workbook automatically.
$ws.Range("C3").Value = "RB005"
$ws.Range("C4").Value = "RB005"
$ws.Range("C5").Value = "RB005"

# Move the selection/active cell to F4, matching the cursor position recorded
# in the saved workbook.
[void]$ws.Range("F4").Select()
